$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Preserve the original "Main" sheet as a new sheet "Main Orig"
#    (a full copy of the sheet before any of the new edits are made).
# ------------------------------------------------------------------
$mainSheet = $wb.Worksheets.Item("Main")
$mainSheet.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$origSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$origSheet.Name = "Main Orig"

# The VLOOKUP formula on the copied sheet still refers to the old sheet
# name ("Main"); repoint it at the copy's own name.
$origSheet.Range("J2").Formula = "=VLOOKUP('Main Orig'!J1,Departments!A2:B10,2,FALSE)"

# Mark the last three agents' "Order (by ID)" as excluded ("x") on the
# archived original sheet.
$origSheet.Range("C6").Value = "x"
$origSheet.Range("C7").Value = "x"
$origSheet.Range("C8").Value = "x"

# ------------------------------------------------------------------
# 2) Rework the live "Main" sheet with the new priority-based layout.
# ------------------------------------------------------------------

# Drop the last agent (Shachar Gafni, row 8) entirely.
$mainSheet.Rows("8").Delete()

# Insert a new "Priority On/Off" column before the ticket-count columns.
$mainSheet.Columns("D").Insert()
$mainSheet.Range("D1").Value = "Priority On/Off"

# First three agents are prioritized ("t"), remaining three are not ("f").
$mainSheet.Range("D2").Value = "t"
$mainSheet.Range("D3").Value = "t"
$mainSheet.Range("D4").Value = "t"
$mainSheet.Range("D5").Value = "f"
$mainSheet.Range("D6").Value = "f"
$mainSheet.Range("D7").Value = "f"

# Tickets are now distributed equally (1 each) instead of the old
# hard-coded per-agent maximums.
$mainSheet.Range("E2:G7").Value = 1

# Switch the active department from "Modix" to "Sandbox" (now in column
# K after the column insert shifted I:J to J:K).
$mainSheet.Range("K1").Value = "Sandbox"
$mainSheet.Range("K2").Formula = "=VLOOKUP(Main!K1,Departments!A2:B10,2,FALSE)"
